$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6788.088191908992

# Row 4
$ws.Range("G4").Value = 17086.97919129642

# Row 6
$ws.Range("G6").Value = 22231.83681584428

# Row 8
$ws.Range("G8").Value = 6198.711515966546

# Row 9
$ws.Range("C9").Value = 6
$ws.Range("G9").Value = 5

# Row 10
$ws.Range("G10").Value = 0.4978133868217219

# Row 12
$ws.Range("G12").Value = 0.5797070920391136

# Row 14
$ws.Range("F14").Value = 0.5709397031663825
$ws.Range("G14").Value = 0.5545366382550266

# Row 16
$ws.Range("G16").Value = 0.4927461019866803

# Row 18
$ws.Range("G18").Value = 0.4218885172214436

# Row 20
$ws.Range("G20").Value = 0.4965821404509979

# Row 21
$ws.Range("G21").Value = 6
$ws.Range("I21").Value = 5

# Row 22
$ws.Range("G22").Value = 0.5480300352097685

# Row 24
$ws.Range("G24").Value = 0.5123615864059685

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 4
